$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$rows = @(
    @("lab.vape.button.all.list", "Všechny vapy"),
    @("lab.vape.latest.title", "Nejnovější vapy"),
    @("lab.vape.button.clone", "Klonovat"),
    @("lab.vape.button.index", "Detail vapu"),
    @("lab.vape.clone.title", "Klon vapu"),
    @("lab.build.link.button", "Detail buildu"),
    @("lab.mixture.table.age", "Stáří mixu"),
    @("lab.mixture.steep.done", "Zrání dokončeno")
)

$startRow = 532
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $rows[$i][0]
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]

    # Copy the formatting (wrap-text "import" cell style) from the row
    # above down onto the freshly written row, same as Excel's fill-down.
    $ws.Range("A531:C531").Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3)).PasteSpecial(-4122)
}

$ws.Range("B528").Select()
